$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the last sprint block (Sprint 15, rows 211-222) down to rows 225-236
#    to create the new Sprint 16 block with identical formatting/merges.
$src = $ws.Range("B211:G222")
$dest = $ws.Range("B225")
$src.Copy($dest) | Out-Null

# 2. Fix the label text for the new block.
$ws.Range("B225").Value = "Sprint 16"

# 3. Write the Sprint 16 input data (Min./Real day counts) and formulas.
$ws.Range("D228").Value = 3
$ws.Range("F228").Value = 4
$ws.Range("E228").Formula = "=D228*C228"
$ws.Range("G228").Formula = "=F228*C228"

$ws.Range("D229").Value = 4
$ws.Range("F229").Value = 3
$ws.Range("D230").Value = 5
$ws.Range("F230").Value = 5
$ws.Range("D231").Value = 4
$ws.Range("F231").Value = 4
$ws.Range("D232").Value = 1
$ws.Range("F232").Value = 1
$ws.Range("D233").Value = 0
$ws.Range("F233").Value = 0
$ws.Range("D234").Value = 0
$ws.Range("F234").Value = 0
$ws.Range("D235").Value = 0
$ws.Range("F235").Value = 0

$ws.Range("E229:E235").Formula = "=D229*C229"
$ws.Range("G229:G235").Formula = "=F229*C229"

$ws.Range("E236").Formula = "=SUM(E228:E235)/60"
$ws.Range("G236").Formula = "=SUM(G228:G235)/60"

# 4. Create the Excel Table (ListObject) over the new block.
$tbl = $ws.ListObjects.Add(1, $ws.Range("B227:G236"), $null, 1)
$tbl.Name = "Tabla18101123456712131415161718"
$tbl = $ws.ListObjects.Item("Tabla18101123456712131415161718")
$tbl.TableStyle = "TableStyleMedium16"
$tbl.ShowTableStyleColumnStripes = $true
$tbl.ShowTableStyleRowStripes = $false
$tbl.ListColumns.Item(6).TotalsCalculation = -4157

# 5. Update the grand-total formulas (I3/J3, I4/J4) to include the new sprint.
$ws.Range("J3").Formula = "=E12+E26+E40+E54+E68+E82+E96+E110+E124+E138+E152+E166+E180+E194+E208+E222+E236"
$ws.Range("J4").Formula = "=G12+G40+G54+G68+G82+G96+G110+G124+G138+G152+G166+G180+G194+G208+G222+G236"

# 6. Selection cosmetics
$ws.Range("J4").Select() | Out-Null

Write-Host "Done"
